$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values were converted from fractions to percentages (multiplied by 100)
$ws.Range("D2").Value  = 55.51702395964691
$ws.Range("D3").Value  = 55.0653159157558
$ws.Range("D4").Value  = 55.01931724618405
$ws.Range("D5").Value  = 54.4262738207866
$ws.Range("D6").Value  = 53.53333911384722
$ws.Range("D7").Value  = 52.95366112486735
$ws.Range("D8").Value  = 45.50042052144659
$ws.Range("D9").Value  = 43.46097581963686
$ws.Range("D10").Value = 50.12536671502129
